$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cells (Volume number, week range) ---
$ws.Range("A8").Value = "Volume 29   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/19/2022  Through  12/25/2022"

# --- Crime Complaints table (rows 14-29) ---
# Row 14
$ws.Range("N14").Value = -80
# Row 15
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 200
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 400
$ws.Range("I15").Value = 54
$ws.Range("J15").Value = 54
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -12.903225806451
$ws.Range("M15").Value = -11.475409836065
$ws.Range("N15").Value = -55
# Row 16
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -30.769230769230
$ws.Range("F16").Value = 37
$ws.Range("H16").Value = -27.450980392156
$ws.Range("I16").Value = 740
$ws.Range("J16").Value = 615
$ws.Range("K16").Value = 20.325203252032
$ws.Range("L16").Value = 19.935170178282
$ws.Range("M16").Value = -1.595744680851
$ws.Range("N16").Value = -76.152110860457
# Row 17
$ws.Range("C17").Value = 18
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = -10
$ws.Range("F17").Value = 74
$ws.Range("G17").Value = 78
$ws.Range("H17").Value = -5.128205128205
$ws.Range("I17").Value = 1077
$ws.Range("J17").Value = 949
$ws.Range("K17").Value = 13.487881981032
$ws.Range("L17").Value = 26.260257913247
$ws.Range("M17").Value = 45.540540540540
$ws.Range("N17").Value = -26.384142173615
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 19
$ws.Range("H18").Value = -38.709677419354
$ws.Range("I18").Value = 369
$ws.Range("J18").Value = 399
$ws.Range("K18").Value = -7.518796992481
$ws.Range("L18").Value = -13.380281690140
$ws.Range("M18").Value = -21.656050955414
$ws.Range("N18").Value = -79.780821917808
# Row 19
$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 26
$ws.Range("E19").Value = -11.538461538461
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 112
$ws.Range("H19").Value = -33.928571428571
$ws.Range("I19").Value = 1282
$ws.Range("J19").Value = 1076
$ws.Range("K19").Value = 19.144981412639
$ws.Range("L19").Value = 35.089567966280
$ws.Range("M19").Value = 81.073446327683
$ws.Range("N19").Value = 34.522560335781
# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 125
$ws.Range("F20").Value = 47
$ws.Range("G20").Value = 31
$ws.Range("H20").Value = 51.612903225806
$ws.Range("I20").Value = 515
$ws.Range("J20").Value = 515
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 26.535626535626
$ws.Range("M20").Value = 61.442006269592
$ws.Range("N20").Value = -80.313455657492
# Row 21
$ws.Range("C21").Value = 66
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = -5.714285714285
$ws.Range("F21").Value = 257
$ws.Range("G21").Value = 304
$ws.Range("H21").Value = -15.460526315789
$ws.Range("I21").Value = 4062
$ws.Range("J21").Value = 3633
$ws.Range("K21").Value = 11.808422791081
$ws.Range("L21").Value = 21.616766467065
$ws.Range("M21").Value = 31.797534068786
$ws.Range("N21").Value = -60.195982361587
# Row 22
$ws.Range("C22").Value = 3
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 200
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 77
$ws.Range("J22").Value = 55
$ws.Range("K22").Value = 40
$ws.Range("L22").Value = 18.461538461538
$ws.Range("M22").Value = 2.666666666666
# Row 23
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 80
$ws.Range("F23").Value = 23
$ws.Range("G23").Value = 32
$ws.Range("H23").Value = -28.125
$ws.Range("I23").Value = 349
$ws.Range("J23").Value = 304
$ws.Range("K23").Value = 14.802631578947
$ws.Range("L23").Value = 31.698113207547
$ws.Range("M23").Value = 54.424778761061
# Row 24
$ws.Range("C24").Value = 42
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = -19.230769230769
$ws.Range("F24").Value = 208
$ws.Range("G24").Value = 229
$ws.Range("H24").Value = -9.170305676855
$ws.Range("I24").Value = 2858
$ws.Range("J24").Value = 2344
$ws.Range("K24").Value = 21.928327645051
$ws.Range("L24").Value = 46.039856923863
$ws.Range("M24").Value = 91.298527443105
# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 71
$ws.Range("G25").Value = 101
$ws.Range("H25").Value = -29.702970297029
$ws.Range("I25").Value = 1242
$ws.Range("J25").Value = 1051
$ws.Range("K25").Value = 18.173168411037
$ws.Range("L25").Value = 40.338983050847
$ws.Range("M25").Value = -31.645569620253
# Row 26
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 7
$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 80
$ws.Range("J26").Value = 97
$ws.Range("K26").Value = -17.525773195876
$ws.Range("L26").Value = -10.112359550561
# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = 33.333333333333
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 18
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 115
$ws.Range("J27").Value = 129
$ws.Range("K27").Value = -10.852713178294
$ws.Range("L27").Value = 10.576923076923
# Row 28
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 100
$ws.Range("I28").Value = 90
$ws.Range("J28").Value = 84
$ws.Range("K28").Value = 7.142857142857
$ws.Range("L28").Value = -27.419354838709
$ws.Range("M28").Value = -15.887850467289
$ws.Range("N28").Value = -71.061093247588
# Row 29
$ws.Range("C29").Value = 2
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("E29").Value = 100
$ws.Range("I29").Value = 79
$ws.Range("J29").Value = 70
$ws.Range("K29").Value = 12.857142857142
$ws.Range("L29").Value = -22.549019607843
$ws.Range("M29").Value = -10.227272727272
$ws.Range("N29").Value = -71.985815602836
